# Updates symbol list: refresh prices and re-rank several coin rows (commit: "Updated symbol list").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.22"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.75"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.454"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.432"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8135"
$ws.Range("E7").Value = "6MXTokenMX"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8848"
$ws.Range("E8").Value = "7FTXTokenFTT"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1442"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07352"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("B11").Value = "ProBitToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1320"
$ws.Range("E11").Value = "10ProBitTokenPROB"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02981"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03063"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09395"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001591"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005841"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005108"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009955"
$ws.Range("B23").Value = "KuCoinToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.329"
$ws.Range("E23").Value = "22KuCoinTokenKCS"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.196"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3276"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.165"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003159"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03918"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006758"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007789"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005646"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3801"
